# Update "想去人数" (want-to-go count) figures in the 苏州-漫展信息 workbook.
# Mirrors the upstream data refresh recorded in commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value  = 594
$wsExpo.Range("F7").Value  = 14844
$wsExpo.Range("F9").Value  = 674
$wsExpo.Range("F11").Value = 34
$wsExpo.Range("F12").Value = 8695
$wsExpo.Range("F19").Value = 10
$wsExpo.Range("F21").Value = 148
$wsExpo.Range("F32").Value = 28
$wsExpo.Range("F35").Value = 267
$wsExpo.Range("F38").Value = 5348

# --- Sheet: 演出 (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 1004

# --- Sheet: 全部类型 (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 594
$wsAll.Range("F7").Value  = 14844
$wsAll.Range("F9").Value  = 674
$wsAll.Range("F11").Value = 34
$wsAll.Range("F12").Value = 8695
$wsAll.Range("F15").Value = 1004
$wsAll.Range("F20").Value = 10
$wsAll.Range("F22").Value = 151
$wsAll.Range("F35").Value = 28
$wsAll.Range("F38").Value = 267
$wsAll.Range("F41").Value = 5348
